$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Keywords" worksheet as the last tab.
# ---------------------------------------------------------------------------
$wsModules = $wb.Worksheets.Item("test_modules")
$wsKeywords = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsKeywords.Name = "Keywords"

# Populate A1:A5 with the keyword list (creates shared strings 30-34 in this order).
$wsKeywords.Range("A1").Value = "<td>"
$wsKeywords.Range("A2").Value = "<ti>"
$wsKeywords.Range("A3").Value = "<tl>"
$wsKeywords.Range("A4").Value = "<te>"
$wsKeywords.Range("A5").Value = "<re>"

# Bigger font for the keyword list.
$wsKeywords.Range("A1:A5").Font.Size = 20

# Sort the list alphabetically (this is how the sheet ends up re-ordered as
# <re>, <td>, <te>, <ti>, <tl> while keeping the shared-string table order
# intact) and also leaves a <sortState> behind in the sheet XML.
$wsKeywords.Sort.SortFields.Clear()
$wsKeywords.Sort.SortFields.Add($wsKeywords.Range("A1")) | Out-Null
$wsKeywords.Sort.SetRange($wsKeywords.Range("A1:A5"))
$wsKeywords.Sort.Header = 2
$wsKeywords.Sort.Apply()

$wsKeywords.Range("A1:A5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Fill in the keyword columns on the "test_modules" sheet.
#    Cell values are written in the exact order needed so new shared
#    strings are minted in the same sequence as the target workbook.
# ---------------------------------------------------------------------------

# -- New shared strings 35-53, minted in first-use order --------------------
$wsModules.Range("B2").Value = "none"

$wsModules.Range("N2").Value = "<if>"
$wsModules.Range("N3").Value = "<else if>"
$wsModules.Range("N4").Value = "<else>"
$wsModules.Range("N5").Value = "<end if>"

$wsModules.Range("N7").Value = "<loopexit>"
$wsModules.Range("N8").Value = "<loopend>"
$wsModules.Range("N6").Value = "<loopstart>"

$wsModules.Range("G6").Value = "<click>"
$wsModules.Range("G7").Value = "<doubleclick>"
$wsModules.Range("G8").Value = "<rightclick>"
$wsModules.Range("G9").Value = "<on>"
$wsModules.Range("G10").Value = "<off>"
$wsModules.Range("G11").Value = "<blank>"
$wsModules.Range("G12").Value = "<first>"
$wsModules.Range("G13").Value = "<second>"
$wsModules.Range("G14").Value = "<third>"
$wsModules.Range("G15").Value = "<last>"
$wsModules.Range("G16").Value = "<random>"

# -- Remaining cells, all reusing already-minted shared strings -------------
$wsModules.Range("C2").Value = "<te>"
$wsModules.Range("D2").Value = "none"
$wsModules.Range("E2").Value = "none"
$wsModules.Range("F2").Value = "<re>"
$wsModules.Range("G2").Value = "<re>"
$wsModules.Range("H2").Value = "none"
$wsModules.Range("I2").Value = "none"
$wsModules.Range("J2").Value = "none"
$wsModules.Range("K2").Value = "none"
$wsModules.Range("L2").Value = "none"
$wsModules.Range("M2").Value = "none"
$wsModules.Range("O2").Value = "none"
$wsModules.Range("P2").Value = "none"
$wsModules.Range("Q2").Value = "none"
$wsModules.Range("R2").Value = "none"
$wsModules.Range("S2").Value = "none"

$wsModules.Range("F3").Value = "<td>"
$wsModules.Range("G3").Value = "<td>"
$wsModules.Range("F4").Value = "<ti>"
$wsModules.Range("G4").Value = "<ti>"
$wsModules.Range("F5").Value = "<tl>"
$wsModules.Range("G5").Value = "<tl>"

# ---------------------------------------------------------------------------
# 3. Activate "test_modules" and fix up its selection/pane state.
# ---------------------------------------------------------------------------
$wsModules.Activate()
$wsModules.Range("N1").Select() | Out-Null
